$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$series = $chart.SeriesCollection(1)
$series.Values = "=Sheet1!`$B`$2:`$B`$104"
$series.XValues = "=Sheet1!`$A`$2:`$A`$104"
